$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.898.43"
$ws.Range("E2").Value = "  -0.86%  "
$ws.Range("D3").Value = "1.880.39"
$ws.Range("E3").Value = "  -1.61%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.53%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.60"
$ws.Range("E5").Value = "  -0.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  -0.33%  "
$ws.Range("E7").Value = "  -0.71%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3882"
$ws.Range("E8").Value = "  -1.89%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07857"
$ws.Range("E9").Value = "  -1.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9848"
$ws.Range("E10").Value = "  -2.62%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.78"
$ws.Range("E11").Value = "  -2.11%  "
$ws.Range("D12").Value = "1.856.38"
$ws.Range("E12").Value = "  -3.91%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.006"
$ws.Range("E13").Value = "  -1.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.676"
$ws.Range("E14").Value = "  -2.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06965"
$ws.Range("E15").Value = "  +0.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.68"
$ws.Range("E16").Value = "  -0.68%  "
$ws.Range("E17").Value = "  -0.36%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009936"
$ws.Range("E18").Value = "  -2.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.97"
$ws.Range("E19").Value = "  -1.99%  "
$ws.Range("E20").Value = "  -0.29%  "
$ws.Range("D21").Value = "28.901.17"
$ws.Range("E21").Value = "  -0.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.274"
$ws.Range("E22").Value = "  -2.46%  "
$ws.Range("E23").Value = "  -1.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.104"
$ws.Range("E24").Value = "  +2.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.24"
$ws.Range("E25").Value = "  +0.49%  "
$ws.Range("E26").Value = "  -1.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.907"
$ws.Range("E27").Value = "  -0.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "117.78"
$ws.Range("E28").Value = "  -2.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.903"
$ws.Range("E29").Value = "  -6.24%  "
$ws.Range("E30").Value = "  -0.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.9008"
$ws.Range("E31").Value = "  -4.33%  "
$ws.Range("E32").Value = "  -1.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.321"
$ws.Range("E33").Value = "  -1.88%  "
$ws.Range("E34").Value = "  -0.94%  "
$ws.Range("E35").Value = "  -2.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.172"
$ws.Range("E36").Value = "  +0.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02080"
$ws.Range("E37").Value = "  -1.09%  "
$ws.Range("E38").Value = "  -0.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.642"
$ws.Range("E39").Value = "  -6.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5666"
$ws.Range("E40").Value = "  -3.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1776"
$ws.Range("E41").Value = "  -2.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.694"
$ws.Range("E42").Value = "  -3.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.99"
$ws.Range("E43").Value = "  -0.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.238"
$ws.Range("E44").Value = "  -2.82%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5341"
$ws.Range("E45").Value = "  -2.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07043"
$ws.Range("E46").Value = "  -2.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.846"
$ws.Range("E47").Value = "  -3.46%  "
$ws.Range("E48").Value = "  +1.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "112.45"
$ws.Range("E49").Value = "  -0.51%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.069"
$ws.Range("E50").Value = "  -5.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "70.85"
$ws.Range("E51").Value = "  -1.57%  "
